$wb = $excel.ActiveWorkbook

# The new "Event" column is added to the "Card24" sheet (first sheet in
# the workbook). Select it explicitly in case it isn't already active.
$ws = $wb.Worksheets.Item("Card24")
$ws.Activate()

# New header cell M1 = "Event", formatted like the other header cells
# (bold, centered, bordered) by copying the format from L1 (the current
# last header column) instead of re-declaring a brand new style.
$ws.Range("M1").Value = "Event"
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats

# Extend column M down through the existing data rows (2-12) so the new
# column participates in the sheet's used range, mirroring the formatting
# (i.e. none) of the neighbouring column L cells.
for ($r = 2; $r -le 12; $r++) {
    $srcCell = $ws.Cells.Item($r, 12)   # column L
    $dstCell = $ws.Cells.Item($r, 13)   # column M
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)        # xlPasteFormats
}

$excel.CutCopyMode = 0

Write-Output "Added 'Event' column (M) to Card24"
